$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter newly logged hours for tasks (yellow "entered" cells)
$ws.Range("H9").Value = 1
$ws.Range("H9").Interior.Color = 65535

$ws.Range("L12").Value = 0.5
$ws.Range("L12").Interior.Color = 65535

$ws.Range("L14").Value = 0.5
$ws.Range("L14").Interior.Color = 65535

# L17 ("Actual Estimated Remaining Hours" for Day 10) needs a bespoke
# formula now that more hours have been logged than originally estimated,
# so it no longer follows the simple day-over-day burndown pattern.
$ws.Range("L17").Formula = "=K17-SUM(L3:L14)-(SUM(B3:B14)-SUM(C3:L14))"

# Update the selection to match the post-edit state
$ws.Range("L18").Select()
